$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")

# Remove the two mailto: hyperlinks that lived on N2/N3 (jrua@todo1.net)
$ws.Hyperlinks.Delete()

# Wipe out the "correo" / "tipoCorreo" / "numeroCelular" columns' content.
# N1:O3 held the headers + the two data rows; the cells keep their existing
# formatting (styles s="9","s="10","s="11","s="12") but become blank.
$ws.Range("N1:O3").ClearContents()

# Column P ("numeroCelular" header + the two 3146834995 values) is removed
# entirely, shrinking the used range from A1:P7 down to A1:O7.
$ws.Columns("P:P").Delete()

# numeroDocumento (B2) gets left-aligned explicitly, producing a new cell
# style that layers horizontal="left" on top of the existing border/fill.
$ws.Range("B2").HorizontalAlignment = -4131  # xlLeft

# Update the active selection/view: it used to be frozen scrolled to F1 with
# F2 selected; now nothing is scrolled and D2 is the selected cell.
$ws.Activate()
$ws.Range("D2").Select()
